$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target for this dataset changed from deuteron to proton ("p").
# Column I (rows 2-10) holds the "target" value for each data row.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "p"
}

# Bold the header row (A1:K1), keeping its existing center alignment.
$ws.Range("A1:K1").Font.Bold = $true

# Move the active selection/cursor to H16 (cosmetic, matches author's last position).
$ws.Range("H16").Select()
